$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.398441
$ws.Range("H2").Value = 0.796882
$ws.Range("I2").Value = 0.08945363909080989
$ws.Range("J2").Value = 0.06146862341190577
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.4154945
$ws.Range("N2").Value = 26.830989
$ws.Range("O2").Value = 0.1359286037219744
$ws.Range("P2").Value = 0.1008397935132719
$ws.Range("Q2").Value = 5.3452830440745
$ws.Range("R2").Value = 21.381132176298
$ws.Range("S2").Value = 0.01215930825946321
$ws.Range("T2").Value = 0.006198483292401652

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.398441
$ws.Range("H3").Value = 0.796882
$ws.Range("I3").Value = 0.08945363909080989
$ws.Range("J3").Value = 0.06146862341190577
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.31779433333333
$ws.Range("N3").Value = 60.953383
$ws.Range("O3").Value = 0.2058641531581443
$ws.Range("P3").Value = 0.2290831156337689
$ws.Range("Q3").Value = 8.095442291967666
$ws.Range("R3").Value = 48.57265375180599
$ws.Range("S3").Value = 0.01841529765834385
$ws.Range("T3").Value = 0.0140814237649182

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.398441
$ws.Range("H4").Value = 0.796882
$ws.Range("I4").Value = 0.08945363909080989
$ws.Range("J4").Value = 0.06146862341190577
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.218472
$ws.Range("N4").Value = 54.655416
$ws.Range("O4").Value = 0.1845933790146823
$ws.Range("P4").Value = 0.2054132579243344
$ws.Range("Q4").Value = 7.258986202151999
$ws.Range("R4").Value = 43.553917212912
$ws.Range("S4").Value = 0.01651254950493247
$ws.Range("T4").Value = 0.01262647019516358

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.398441
$ws.Range("H5").Value = 0.796882
$ws.Range("I5").Value = 0.08945363909080989
$ws.Range("J5").Value = 0.06146862341190577
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.09656133333333
$ws.Range("N5").Value = 60.289684
$ws.Range("O5").Value = 0.2036225740059764
$ws.Range("P5").Value = 0.2265887137272657
$ws.Range("Q5").Value = 8.007293994214667
$ws.Range("R5").Value = 48.043763965288
$ws.Range("S5").Value = 0.01821478024587234
$ws.Range("T5").Value = 0.01392809631348942

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.398441
$ws.Range("H6").Value = 0.796882
$ws.Range("I6").Value = 0.08945363909080989
$ws.Range("J6").Value = 0.06146862341190577
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.05227
$ws.Range("N6").Value = 30.15681
$ws.Range("O6").Value = 0.1018517077649498
$ws.Range("P6").Value = 0.1133393365939278
$ws.Range("Q6").Value = 4.00523651107
$ws.Range("R6").Value = 24.03141906642
$ws.Range("S6").Value = 0.009111005907188456
$ws.Range("T6").Value = 0.006966812998847381

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.398441
$ws.Range("H7").Value = 0.796882
$ws.Range("I7").Value = 0.08945363909080989
$ws.Range("J7").Value = 0.06146862341190577
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.594562
$ws.Range("N7").Value = 33.189124
$ws.Range("O7").Value = 0.1681395823342728
$ws.Range("P7").Value = 0.1247357826074312
$ws.Range("Q7").Value = 6.611953877842
$ws.Range("R7").Value = 26.447815511368
$ws.Range("S7").Value = 0.01504069751500955
$ws.Range("T7").Value = 0.007667336847085536

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.055721
$ws.Range("H8").Value = 12.167163
$ws.Range("I8").Value = 0.9105463609091901
$ws.Range("J8").Value = 0.9385313765880943
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.4154945
$ws.Range("N8").Value = 26.830989
$ws.Range("O8").Value = 0.1359286037219744
$ws.Range("P8").Value = 0.1008397935132719
$ws.Range("Q8").Value = 54.4095027690345
$ws.Range("R8").Value = 326.457016614207
$ws.Range("S8").Value = 0.1237692954625111
$ws.Range("T8").Value = 0.0946413102208703

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.055721
$ws.Range("H9").Value = 12.167163
$ws.Range("I9").Value = 0.9105463609091901
$ws.Range("J9").Value = 0.9385313765880943
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 20.31779433333333
$ws.Range("N9").Value = 60.953383
$ws.Range("O9").Value = 0.2058641531581443
$ws.Range("P9").Value = 0.2290831156337689
$ws.Range("Q9").Value = 82.40330515138099
$ws.Range("R9").Value = 741.6297463624289
$ws.Range("S9").Value = 0.1874488554998004
$ws.Range("T9").Value = 0.2150016918688507

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.055721
$ws.Range("H10").Value = 12.167163
$ws.Range("I10").Value = 0.9105463609091901
$ws.Range("J10").Value = 0.9385313765880943
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 18.218472
$ws.Range("N10").Value = 54.655416
$ws.Range("O10").Value = 0.1845933790146823
$ws.Range("P10").Value = 0.2054132579243344
$ws.Range("Q10").Value = 73.889039478312
$ws.Range("R10").Value = 665.0013553048079
$ws.Range("S10").Value = 0.1680808295097498
$ws.Range("T10").Value = 0.1927867877291708

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.055721
$ws.Range("H11").Value = 12.167163
$ws.Range("I11").Value = 0.9105463609091901
$ws.Range("J11").Value = 0.9385313765880943
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 20.09656133333333
$ws.Range("N11").Value = 60.289684
$ws.Range("O11").Value = 0.2036225740059764
$ws.Range("P11").Value = 0.2265887137272657
$ws.Range("Q11").Value = 81.50604582738801
$ws.Range("R11").Value = 733.554412446492
$ws.Range("S11").Value = 0.185407793760104
$ws.Range("T11").Value = 0.2126606174137763

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.055721
$ws.Range("H12").Value = 12.167163
$ws.Range("I12").Value = 0.9105463609091901
$ws.Range("J12").Value = 0.9385313765880943
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 10.05227
$ws.Range("N12").Value = 30.15681
$ws.Range("O12").Value = 0.1018517077649498
$ws.Range("P12").Value = 0.1133393365939278
$ws.Range("Q12").Value = 40.76920253667
$ws.Range("R12").Value = 366.92282283003
$ws.Range("S12").Value = 0.09274070185776132
$ws.Range("T12").Value = 0.1063725235950804

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.055721
$ws.Range("H13").Value = 12.167163
$ws.Range("I13").Value = 0.9105463609091901
$ws.Range("J13").Value = 0.9385313765880943
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 16.594562
$ws.Range("N13").Value = 33.189124
$ws.Range("O13").Value = 0.1681395823342728
$ws.Range("P13").Value = 0.1247357826074312
$ws.Range("Q13").Value = 67.302913589202
$ws.Range("R13").Value = 403.817481535212
$ws.Range("S13").Value = 0.1530988848192633
$ws.Range("T13").Value = 0.1170684457603457
